# Reorders the player rows (rows 4-16) on the active sheet to match the
# updated roster ordering from the diff. The set of rows/values is the
# same, only the row order changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam", "SF,PF", "Indiana Pacers"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jakob Poeltl", "C", "Toronto Raptors"),
    @("Clint Capela", "C", "Atlanta Hawks"),
    @("Al Horford", "PF,C", "Boston Celtics"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Russell Westbrook", "PG", "Denver Nuggets"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks")
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
